$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1164814.6
$ws.Range("J17").Value = 1164814.6
$ws.Range("L17").Value = 3494443.8
$ws.Range("N17").Value = -3494779.8
$ws.Range("H21").Value = 12003.4
$ws.Range("J21").Value = 27500
$ws.Range("L21").Value = 27500
$ws.Range("N21").Value = -28436
$ws.Range("H23").Value = 12003.4
$ws.Range("J23").Value = 27500
$ws.Range("L23").Value = 27500
$ws.Range("N23").Value = -27968
$ws.Range("H112").Value = 22728674
$ws.Range("J112").Value = 30304568
$ws.Range("L112").Value = 90913704
$ws.Range("N112").Value = -90915920
$ws.Range("H129").Value = 955.51514
$ws.Range("J129").Value = 1105.6923
$ws.Range("L129").Value = 3317.0769
$ws.Range("N129").Value = -13317.0769
$ws.Range("H138").Value = 6889990
$ws.Range("I138").Value = 5056037
$ws.Range("J138").Value = 7248807
$ws.Range("K138").Value = 15168111
$ws.Range("L138").Value = 21746421
$ws.Range("M138").Value = -15162971
$ws.Range("N138").Value = -21756701

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21473.781
$ws.Range("I32").Value = 2804.9805
$ws.Range("K32").Value = 2804.9805
$ws.Range("M32").Value = -2517.9805
$ws.Range("H45").Value = 1870.0834
$ws.Range("I45").Value = 1033
$ws.Range("K45").Value = 1033
$ws.Range("M45").Value = -656
$ws.Range("H61").Value = 2507.5789
$ws.Range("I61").Value = 1759.7587
$ws.Range("J61").Value = 4917.222
$ws.Range("K61").Value = 1759.7587
$ws.Range("L61").Value = 4917.222
$ws.Range("M61").Value = -1547.7587
$ws.Range("N61").Value = -5341.222
$ws.Range("H109").Value = 30000
$ws.Range("I109").Value = 30000
$ws.Range("K109").Value = 30000
$ws.Range("M109").Value = -28613
$ws.Range("H122").Value = 1503.5526
$ws.Range("I122").Value = 1437.8966
$ws.Range("J122").Value = 1715.1111
$ws.Range("K122").Value = 4313.6898
$ws.Range("L122").Value = 5145.3333
$ws.Range("M122").Value = -1863.6898
$ws.Range("N122").Value = -10045.3333
$ws.Range("H132").Value = 2649.625
$ws.Range("I132").Value = 2249.9348
$ws.Range("J132").Value = 4488.2
$ws.Range("K132").Value = 6749.8044
$ws.Range("L132").Value = 13464.6
$ws.Range("M132").Value = -4219.8044
$ws.Range("N132").Value = -18524.6
$ws.Range("H136").Value = 2507.5789
$ws.Range("I136").Value = 1759.7587
$ws.Range("J136").Value = 4917.222
$ws.Range("K136").Value = 5279.2761
$ws.Range("L136").Value = 14751.666
$ws.Range("M136").Value = -2729.2761
$ws.Range("N136").Value = -19851.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 44000
$ws.Range("J58").Value = 44000
$ws.Range("L58").Value = 44000
$ws.Range("N58").Value = -44588
$ws.Range("H134").Value = 25643470
$ws.Range("I134").Value = 33334676
$ws.Range("J134").Value = 6118
$ws.Range("K134").Value = 100004028
$ws.Range("L134").Value = 18354
$ws.Range("M134").Value = -100001493
$ws.Range("N134").Value = -23424
$ws.Range("H135").Value = 50000
$ws.Range("J135").Value = 50000
$ws.Range("L135").Value = 50000
$ws.Range("N135").Value = -60140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 370.42856
$ws.Range("I22").Value = 382.16666
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 382.16666
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = -32.16665999999998
$ws.Range("N22").Value = -1000
$ws.Range("H58").Value = 2304.1667
$ws.Range("I58").Value = 1637.6666
$ws.Range("J58").Value = 4303.6665
$ws.Range("K58").Value = 1637.6666
$ws.Range("L58").Value = 4303.6665
$ws.Range("M58").Value = -1434.6666
$ws.Range("N58").Value = -4709.6665
$ws.Range("H132").Value = 2896.3428
$ws.Range("I132").Value = 2470.7585
$ws.Range("K132").Value = 7412.2755
$ws.Range("M132").Value = -4882.2755
$ws.Range("H134").Value = 3342.037
$ws.Range("I134").Value = 1714.2941
$ws.Range("J134").Value = 6109.2
$ws.Range("K134").Value = 5142.8823
$ws.Range("L134").Value = 18327.6
$ws.Range("M134").Value = -2607.8823
$ws.Range("N134").Value = -23397.6
$ws.Range("H136").Value = 2304.1667
$ws.Range("I136").Value = 1637.6666
$ws.Range("J136").Value = 4303.6665
$ws.Range("K136").Value = 4912.9998
$ws.Range("L136").Value = 12910.9995
$ws.Range("M136").Value = -2362.9998
$ws.Range("N136").Value = -18010.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 1814.5
$ws.Range("I116").Value = 629
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 1887
$ws.Range("L116").Value = 9000
$ws.Range("M116").Value = 1555
$ws.Range("N116").Value = -15884
$ws.Range("H129").Value = 1271.6666
$ws.Range("I129").Value = 526
$ws.Range("J129").Value = 5000
$ws.Range("K129").Value = 1578
$ws.Range("L129").Value = 15000
$ws.Range("M129").Value = 3422
$ws.Range("N129").Value = -25000
$ws.Range("H131").Value = 6537366.5
$ws.Range("J131").Value = 6668107.5
$ws.Range("L131").Value = 20004322.5
$ws.Range("N131").Value = -20014402.5
$ws.Range("H132").Value = 8772829
$ws.Range("I132").Value = 579.8
$ws.Range("J132").Value = 11905775
$ws.Range("K132").Value = 5218.2
$ws.Range("L132").Value = 107151975
$ws.Range("M132").Value = -2688.2
$ws.Range("N132").Value = -107157035

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3583333.8
$ws.Range("I11").Value = 3583333.8
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 3583333.8
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -3583194.8
$ws.Range("N11").ClearContents()
$ws.Range("H18").Value = 7100
$ws.Range("J18").Value = 7100
$ws.Range("L18").Value = 7100
$ws.Range("N18").Value = -7686
$ws.Range("H80").Value = 2619.5217
$ws.Range("I80").Value = 2466.611
$ws.Range("J80").Value = 3170
$ws.Range("K80").Value = 2466.611
$ws.Range("L80").Value = 3170
$ws.Range("M80").Value = -1468.611
$ws.Range("N80").Value = -5166
$ws.Range("H83").Value = 2619.5217
$ws.Range("I83").Value = 2466.611
$ws.Range("J83").Value = 3170
$ws.Range("K83").Value = 12333.055
$ws.Range("L83").Value = 15850
$ws.Range("M83").Value = -7341.055
$ws.Range("N83").Value = -25834
$ws.Range("H93").Value = 26848.857
$ws.Range("J93").Value = 26848.857
$ws.Range("L93").Value = 26848.857
$ws.Range("N93").Value = -30592.857
$ws.Range("H109").Value = 132785
$ws.Range("J109").Value = 132785
$ws.Range("L109").Value = 132785
$ws.Range("N109").Value = -134865
$ws.Range("H122").Value = 587090.25
$ws.Range("I122").Value = 1112908.9
$ws.Range("J122").Value = 2847.3333
$ws.Range("K122").Value = 3338726.7
$ws.Range("L122").Value = 8541.999899999999
$ws.Range("M122").Value = -3336276.7
$ws.Range("N122").Value = -13441.9999
$ws.Range("H123").Value = 9880.421
$ws.Range("J123").Value = 9880.421
$ws.Range("L123").Value = 9880.421
$ws.Range("N123").Value = -14780.421

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 42000
$ws.Range("I23").Value = 42500
$ws.Range("J23").Value = 40000
$ws.Range("K23").Value = 42500
$ws.Range("L23").Value = 40000
$ws.Range("M23").Value = -42270
$ws.Range("N23").Value = -40460
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H25").Value = 33247.5
$ws.Range("I25").Value = 33247.5
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 33247.5
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -33017.5
$ws.Range("N25").ClearContents()
$ws.Range("H58").Value = 10000
$ws.Range("I58").Value = 10000
$ws.Range("K58").Value = 10000
$ws.Range("M58").Value = -9740
$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 1000
$ws.Range("K107").Value = 1000
$ws.Range("M107").Value = 920
$ws.Range("H110").Value = 36585.8
$ws.Range("J110").Value = 36585.8
$ws.Range("L110").Value = 36585.8
$ws.Range("N110").Value = -44765.8
$ws.Range("H122").Value = 3007.5
$ws.Range("I122").Value = 2053.4614
$ws.Range("K122").Value = 6160.3842
$ws.Range("M122").Value = -3710.3842

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H53").Value = 11660
$ws.Range("J53").Value = 7400
$ws.Range("L53").Value = 7400
$ws.Range("N53").Value = -8614
$ws.Range("H122").Value = 167834
$ws.Range("I122").Value = 250751
$ws.Range("K122").Value = 752253
$ws.Range("M122").Value = -749803
$ws.Range("H126").Value = 111675.78
$ws.Range("I126").Value = 111675.78
$ws.Range("K126").Value = 335027.34
$ws.Range("M126").Value = -332557.34
$ws.Range("H132").Value = 2827.75
$ws.Range("I132").Value = 2803.875
$ws.Range("J132").Value = 2851.625
$ws.Range("K132").Value = 8411.625
$ws.Range("L132").Value = 8554.875
$ws.Range("M132").Value = -5881.625
$ws.Range("N132").Value = -13614.875
$ws.Range("H135").Value = 143928.75
$ws.Range("J135").Value = 143928.75
$ws.Range("L135").Value = 143928.75
$ws.Range("N135").Value = -154068.75
